$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows from bottom to top so earlier row numbers stay valid while deleting.
$ws.Rows.Item(14).Delete()   # 000330949 Renato      8000
$ws.Rows.Item(10).Delete()   # 004342617 Juraci       14216.8
$ws.Rows.Item(9).Delete()    # 004221638 Caroline     14301.44
$ws.Rows.Item(7).Delete()    # 004211922 Carlos       18627.98
$ws.Rows.Item(6).Delete()    # 004477812 Diego        21018.21
$ws.Rows.Item(3).Delete()    # 004207955 Silvania     79757.93

# After the six deletions above, the row that used to hold
# 005993550 / Alessandra / 1454.14 (originally row 17) has shifted up to row 11.
# Replace its contents with the new record.
$ws.Cells.Item(11, 1).NumberFormat = "@"
$ws.Cells.Item(11, 1).Value = "005135281"
$ws.Cells.Item(11, 2).Value = "Rafael"
$ws.Cells.Item(11, 3).Value = 1500
